$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Cells.Item(70, 8).Value = 10000  # H70: 0 -> 10000
$ws.Cells.Item(70, 9).Value = 10000  # I70: 0 -> 10000
$ws.Cells.Item(70, 11).Value = 30000  # K70: 0 -> 30000
$ws.Cells.Item(70, 13).Value = -29730  # M70: None -> -29730
# Row 73
$ws.Cells.Item(73, 8).Value = 10000  # H73: 0 -> 10000
$ws.Cells.Item(73, 9).Value = 10000  # I73: 0 -> 10000
$ws.Cells.Item(73, 11).Value = 30000  # K73: 0 -> 30000
$ws.Cells.Item(73, 13).Value = -29064  # M73: None -> -29064
# Row 103
$ws.Cells.Item(103, 8).Value = 83334320  # H103: 71429850 -> 83334320
$ws.Cells.Item(103, 10).Value = 2000  # J103: 2500 -> 2000
$ws.Cells.Item(103, 12).Value = 6000  # L103: 7500 -> 6000
$ws.Cells.Item(103, 14).Value = -7172  # N103: -8672 -> -7172
# Row 132
$ws.Cells.Item(132, 8).Value = 5387.8335  # H132: 5183.0527 -> 5387.8335
$ws.Cells.Item(132, 9).Value = 5436.5625  # I132: 5204.8237 -> 5436.5625
$ws.Cells.Item(132, 11).Value = 16309.6875  # K132: 15614.4711 -> 16309.6875
$ws.Cells.Item(132, 13).Value = -13779.6875  # M132: -13084.4711 -> -13779.6875
# Row 135
$ws.Cells.Item(135, 8).Value = 2298  # H135: 2015.3334 -> 2298
$ws.Cells.Item(135, 9).Value = 499.5  # I135: 524.75 -> 499.5
$ws.Cells.Item(135, 10).Value = 3497  # J135: 4996.5 -> 3497
$ws.Cells.Item(135, 11).Value = 4495.5  # K135: 4722.75 -> 4495.5
$ws.Cells.Item(135, 12).Value = 31473  # L135: 44968.5 -> 31473
$ws.Cells.Item(135, 13).Value = -1960.5  # M135: -2187.75 -> -1960.5
$ws.Cells.Item(135, 14).Value = -36543  # N135: -50038.5 -> -36543
# Row 141
$ws.Cells.Item(141, 8).Value = 18047.5  # H141: 9447.5 -> 18047.5
$ws.Cells.Item(141, 9).Value = 0  # I141: 800 -> 0
$ws.Cells.Item(141, 10).Value = 18047.5  # J141: 18095 -> 18047.5
$ws.Cells.Item(141, 11).Value = 0  # K141: 2400 -> 0
$ws.Cells.Item(141, 12).Value = 54142.5  # L141: 54285 -> 54142.5
$ws.Cells.Item(141, 13).ClearContents()  # M141: was 2780
$ws.Cells.Item(141, 14).Value = -64502.5  # N141: -64645 -> -64502.5

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Cells.Item(5, 8).Value = 569.8333  # H5: 643.8 -> 569.8333
$ws.Cells.Item(5, 9).Value = 483.8  # I5: 554.75 -> 483.8
$ws.Cells.Item(5, 11).Value = 483.8  # K5: 554.75 -> 483.8
$ws.Cells.Item(5, 13).Value = -371.8  # M5: -442.75 -> -371.8
# Row 45
$ws.Cells.Item(45, 8).Value = 1174.875  # H45: 1109.3 -> 1174.875
$ws.Cells.Item(45, 9).Value = 1049.8334  # I45: 999.125 -> 1049.8334
$ws.Cells.Item(45, 11).Value = 1049.8334  # K45: 999.125 -> 1049.8334
$ws.Cells.Item(45, 13).Value = -672.8334  # M45: -622.125 -> -672.8334
# Row 61
$ws.Cells.Item(61, 8).Value = 4643.4116  # H61: 4364.4736 -> 4643.4116
$ws.Cells.Item(61, 9).Value = 4643.4116  # I61: 4364.4736 -> 4643.4116
$ws.Cells.Item(61, 11).Value = 4643.4116  # K61: 4364.4736 -> 4643.4116
$ws.Cells.Item(61, 13).Value = -4431.4116  # M61: -4152.4736 -> -4431.4116
# Row 110
$ws.Cells.Item(110, 8).Value = 8501.5  # H110: 6569.625 -> 8501.5
$ws.Cells.Item(110, 9).Value = 14336.667  # I110: 8911.6 -> 14336.667
$ws.Cells.Item(110, 11).Value = 14336.667  # K110: 8911.6 -> 14336.667
$ws.Cells.Item(110, 13).Value = -12291.667  # M110: -6866.6 -> -12291.667
# Row 122
$ws.Cells.Item(122, 8).Value = 1504  # H122: 2006 -> 1504
$ws.Cells.Item(122, 9).Value = 1504  # I122: 2006 -> 1504
$ws.Cells.Item(122, 11).Value = 4512  # K122: 6018 -> 4512
$ws.Cells.Item(122, 13).Value = -2062  # M122: -3568 -> -2062
# Row 132
$ws.Cells.Item(132, 8).Value = 1829.3793  # H132: 1917.875 -> 1829.3793
$ws.Cells.Item(132, 9).Value = 1813  # I132: 1896.4445 -> 1813
$ws.Cells.Item(132, 10).Value = 1865.7778  # J132: 1982.1666 -> 1865.7778
$ws.Cells.Item(132, 11).Value = 5439  # K132: 5689.333500000001 -> 5439
$ws.Cells.Item(132, 12).Value = 5597.3334  # L132: 5946.4998 -> 5597.3334
$ws.Cells.Item(132, 13).Value = -2909  # M132: -3159.333500000001 -> -2909
$ws.Cells.Item(132, 14).Value = -10657.3334  # N132: -11006.4998 -> -10657.3334
# Row 136
$ws.Cells.Item(136, 8).Value = 4643.4116  # H136: 4364.4736 -> 4643.4116
$ws.Cells.Item(136, 9).Value = 4643.4116  # I136: 4364.4736 -> 4643.4116
$ws.Cells.Item(136, 11).Value = 13930.2348  # K136: 13093.4208 -> 13930.2348
$ws.Cells.Item(136, 13).Value = -11380.2348  # M136: -10543.4208 -> -11380.2348

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 25
$ws.Cells.Item(25, 8).Value = 6745  # H25: 7195 -> 6745
$ws.Cells.Item(25, 9).Value = 4993.3335  # I25: 5990 -> 4993.3335
$ws.Cells.Item(25, 10).Value = 12000  # J25: 8400 -> 12000
$ws.Cells.Item(25, 11).Value = 4993.3335  # K25: 5990 -> 4993.3335
$ws.Cells.Item(25, 12).Value = 12000  # L25: 8400 -> 12000
$ws.Cells.Item(25, 13).Value = -4758.3335  # M25: -5755 -> -4758.3335
$ws.Cells.Item(25, 14).Value = -12470  # N25: -8870 -> -12470
# Row 55
$ws.Cells.Item(55, 8).Value = 40000  # H55: 0 -> 40000
$ws.Cells.Item(55, 10).Value = 40000  # J55: 0 -> 40000
$ws.Cells.Item(55, 12).Value = 40000  # L55: 0 -> 40000
$ws.Cells.Item(55, 14).Value = -40546  # N55: None -> -40546
# Row 82
$ws.Cells.Item(82, 8).Value = 6892.857  # H82: 6893.143 -> 6892.857
$ws.Cells.Item(82, 9).Value = 6892.857  # I82: 6893.143 -> 6892.857
$ws.Cells.Item(82, 11).Value = 6892.857  # K82: 6893.143 -> 6892.857
$ws.Cells.Item(82, 13).Value = -6509.857  # M82: -6510.143 -> -6509.857
# Row 85
$ws.Cells.Item(85, 8).Value = 6892.857  # H85: 6893.143 -> 6892.857
$ws.Cells.Item(85, 9).Value = 6892.857  # I85: 6893.143 -> 6892.857
$ws.Cells.Item(85, 11).Value = 6892.857  # K85: 6893.143 -> 6892.857
$ws.Cells.Item(85, 13).Value = -5566.857  # M85: -5567.143 -> -5566.857
# Row 99
$ws.Cells.Item(99, 8).Value = 27780088  # H99: 27780110 -> 27780088
$ws.Cells.Item(99, 9).Value = 8548704  # I99: 9260980 -> 8548704
$ws.Cells.Item(99, 10).Value = 111116080  # J99: 83337496 -> 111116080
$ws.Cells.Item(99, 11).Value = 8548704  # K99: 9260980 -> 8548704
$ws.Cells.Item(99, 12).Value = 111116080  # L99: 83337496 -> 111116080
$ws.Cells.Item(99, 13).Value = -8547206  # M99: -9259482 -> -8547206
$ws.Cells.Item(99, 14).Value = -111119076  # N99: -83340492 -> -111119076
# Row 134
$ws.Cells.Item(134, 8).Value = 5530.4287  # H134: 5439.125 -> 5530.4287
$ws.Cells.Item(134, 9).Value = 4785.5  # I134: 4787.5713 -> 4785.5
$ws.Cells.Item(134, 11).Value = 14356.5  # K134: 14362.7139 -> 14356.5
$ws.Cells.Item(134, 13).Value = -11821.5  # M134: -11827.7139 -> -11821.5
# Row 135
$ws.Cells.Item(135, 8).Value = 0  # H135: 89999 -> 0
$ws.Cells.Item(135, 10).Value = 0  # J135: 89999 -> 0
$ws.Cells.Item(135, 12).Value = 0  # L135: 89999 -> 0
$ws.Cells.Item(135, 14).ClearContents()  # N135: was -100139

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 5359.4  # H31: 5449.25 -> 5359.4
$ws.Cells.Item(31, 10).Value = 6500  # J31: 8000 -> 6500
$ws.Cells.Item(31, 12).Value = 6500  # L31: 8000 -> 6500
$ws.Cells.Item(31, 14).Value = -7090  # N31: -8590 -> -7090
# Row 34
$ws.Cells.Item(34, 8).Value = 5359.4  # H34: 5449.25 -> 5359.4
$ws.Cells.Item(34, 10).Value = 6500  # J34: 8000 -> 6500
$ws.Cells.Item(34, 12).Value = 6500  # L34: 8000 -> 6500
$ws.Cells.Item(34, 14).Value = -6904  # N34: -8404 -> -6904
# Row 107
$ws.Cells.Item(107, 8).Value = 1946.7727  # H107: 1991.9445 -> 1946.7727
$ws.Cells.Item(107, 9).Value = 1866.25  # I107: 1904.6154 -> 1866.25
$ws.Cells.Item(107, 10).Value = 2161.5  # J107: 2219 -> 2161.5
$ws.Cells.Item(107, 11).Value = 1866.25  # K107: 1904.6154 -> 1866.25
$ws.Cells.Item(107, 12).Value = 2161.5  # L107: 2219 -> 2161.5
$ws.Cells.Item(107, 13).Value = 53.75  # M107: 15.38460000000009 -> 53.75
$ws.Cells.Item(107, 14).Value = -6001.5  # N107: -6059 -> -6001.5
# Row 134
$ws.Cells.Item(134, 8).Value = 2951  # H134: 3717 -> 2951
$ws.Cells.Item(134, 9).Value = 3095.9167  # I134: 3717 -> 3095.9167
$ws.Cells.Item(134, 10).Value = 1212  # J134: 0 -> 1212
$ws.Cells.Item(134, 11).Value = 9287.750100000001  # K134: 11151 -> 9287.750100000001
$ws.Cells.Item(134, 12).Value = 3636  # L134: 0 -> 3636
$ws.Cells.Item(134, 13).Value = -6752.750100000001  # M134: -8616 -> -6752.750100000001
$ws.Cells.Item(134, 14).Value = -8706  # N134: None -> -8706

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Cells.Item(2, 8).Value = 170.2  # H2: 188.23077 -> 170.2
$ws.Cells.Item(2, 9).Value = 217.2  # I2: 258.25 -> 217.2
$ws.Cells.Item(2, 11).Value = 217.2  # K2: 258.25 -> 217.2
$ws.Cells.Item(2, 13).Value = -104.2  # M2: -145.25 -> -104.2
# Row 102
$ws.Cells.Item(102, 8).Value = 2594.1  # H102: 2847.7778 -> 2594.1
$ws.Cells.Item(102, 9).Value = 2660.111  # I102: 2953.75 -> 2660.111
$ws.Cells.Item(102, 11).Value = 2660.111  # K102: 2953.75 -> 2660.111
$ws.Cells.Item(102, 13).Value = -1038.111  # M102: -1331.75 -> -1038.111
# Row 122
$ws.Cells.Item(122, 8).Value = 3240.9092  # H122: 3220.4167 -> 3240.9092
$ws.Cells.Item(122, 10).Value = 5633.3335  # J122: 4973.75 -> 5633.3335
$ws.Cells.Item(122, 12).Value = 16900.0005  # L122: 14921.25 -> 16900.0005
$ws.Cells.Item(122, 14).Value = -21800.0005  # N122: -19821.25 -> -21800.0005
# Row 132
$ws.Cells.Item(132, 8).Value = 3186.8333  # H132: 3353.8333 -> 3186.8333
$ws.Cells.Item(132, 9).Value = 3022.75  # I132: 3210.625 -> 3022.75
$ws.Cells.Item(132, 11).Value = 9068.25  # K132: 9631.875 -> 9068.25
$ws.Cells.Item(132, 13).Value = -6538.25  # M132: -7101.875 -> -6538.25

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 11
$ws.Cells.Item(11, 8).Value = 10000  # H11: 1302 -> 10000
$ws.Cells.Item(11, 10).Value = 10000  # J11: 1302 -> 10000
$ws.Cells.Item(11, 12).Value = 10000  # L11: 1302 -> 10000
$ws.Cells.Item(11, 14).Value = -10280  # N11: -1582 -> -10280
# Row 40
$ws.Cells.Item(40, 8).Value = 9283.929  # H40: 5961.75 -> 9283.929
$ws.Cells.Item(40, 9).Value = 6189.5835  # I40: 5961.75 -> 6189.5835
$ws.Cells.Item(40, 10).Value = 27850  # J40: 0 -> 27850
$ws.Cells.Item(40, 11).Value = 6189.5835  # K40: 5961.75 -> 6189.5835
$ws.Cells.Item(40, 12).Value = 27850  # L40: 0 -> 27850
$ws.Cells.Item(40, 13).Value = -6053.5835  # M40: -5825.75 -> -6053.5835
$ws.Cells.Item(40, 14).Value = -28122  # N40: None -> -28122
# Row 46
$ws.Cells.Item(46, 8).Value = 1899  # H46: 1661.6666 -> 1899
$ws.Cells.Item(46, 9).Value = 1899  # I46: 1661.6666 -> 1899
$ws.Cells.Item(46, 11).Value = 1899  # K46: 1661.6666 -> 1899
$ws.Cells.Item(46, 13).Value = -1711  # M46: -1473.6666 -> -1711
# Row 132
$ws.Cells.Item(132, 8).Value = 3203.2  # H132: 3501.5 -> 3203.2
$ws.Cells.Item(132, 9).Value = 3252.75  # I132: 3668.6667 -> 3252.75
$ws.Cells.Item(132, 10).Value = 3005  # J132: 3000 -> 3005
$ws.Cells.Item(132, 11).Value = 9758.25  # K132: 11006.0001 -> 9758.25
$ws.Cells.Item(132, 12).Value = 9015  # L132: 9000 -> 9015
$ws.Cells.Item(132, 13).Value = -7228.25  # M132: -8476.000100000001 -> -7228.25
$ws.Cells.Item(132, 14).Value = -14075  # N132: -14060 -> -14075
# Row 136
$ws.Cells.Item(136, 8).Value = 100004536  # H136: 83337370 -> 100004536
$ws.Cells.Item(136, 9).Value = 9350  # I136: 6733.3335 -> 9350
$ws.Cells.Item(136, 11).Value = 28050  # K136: 20200.0005 -> 28050
$ws.Cells.Item(136, 13).Value = -25500  # M136: -17650.0005 -> -25500
# Row 138
$ws.Cells.Item(138, 8).Value = 50000  # H138: 0 -> 50000
$ws.Cells.Item(138, 10).Value = 50000  # J138: 0 -> 50000
$ws.Cells.Item(138, 12).Value = 50000  # L138: 0 -> 50000
$ws.Cells.Item(138, 14).Value = -60280  # N138: None -> -60280

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 7
$ws.Cells.Item(7, 8).Value = 4500  # H7: 2204 -> 4500
$ws.Cells.Item(7, 9).Value = 4500  # I7: 2204 -> 4500
$ws.Cells.Item(7, 11).Value = 4500  # K7: 2204 -> 4500
$ws.Cells.Item(7, 13).Value = -4387  # M7: -2091 -> -4387
# Row 23
$ws.Cells.Item(23, 8).Value = 3007  # H23: 4103.8 -> 3007
$ws.Cells.Item(23, 9).Value = 3604.75  # I23: 4504.75 -> 3604.75
$ws.Cells.Item(23, 10).Value = 616  # J23: 2500 -> 616
$ws.Cells.Item(23, 11).Value = 3604.75  # K23: 4504.75 -> 3604.75
$ws.Cells.Item(23, 12).Value = 616  # L23: 2500 -> 616
$ws.Cells.Item(23, 13).Value = -3375.75  # M23: -4275.75 -> -3375.75
$ws.Cells.Item(23, 14).Value = -1074  # N23: -2958 -> -1074
# Row 25
$ws.Cells.Item(25, 8).Value = 0  # H25: 25000 -> 0
$ws.Cells.Item(25, 10).Value = 0  # J25: 25000 -> 0
$ws.Cells.Item(25, 12).Value = 0  # L25: 25000 -> 0
$ws.Cells.Item(25, 14).ClearContents()  # N25: was -25586
# Row 81
$ws.Cells.Item(81, 8).Value = 2471.5625  # H81: 2571.0667 -> 2471.5625
$ws.Cells.Item(81, 9).Value = 2110.5715  # I81: 2197.6155 -> 2110.5715
$ws.Cells.Item(81, 11).Value = 4221.143  # K81: 4395.231 -> 4221.143
$ws.Cells.Item(81, 13).Value = -3160.143  # M81: -3334.231 -> -3160.143
# Row 82
$ws.Cells.Item(82, 8).Value = 0  # H82: 32000 -> 0
$ws.Cells.Item(82, 10).Value = 0  # J82: 32000 -> 0
$ws.Cells.Item(82, 12).Value = 0  # L82: 32000 -> 0
$ws.Cells.Item(82, 14).ClearContents()  # N82: was -32766
# Row 84
$ws.Cells.Item(84, 8).Value = 2471.5625  # H84: 2571.0667 -> 2471.5625
$ws.Cells.Item(84, 9).Value = 2110.5715  # I84: 2197.6155 -> 2110.5715
$ws.Cells.Item(84, 11).Value = 21105.715  # K84: 21976.155 -> 21105.715
$ws.Cells.Item(84, 13).Value = -15801.715  # M84: -16672.155 -> -15801.715
# Row 85
$ws.Cells.Item(85, 8).Value = 0  # H85: 32000 -> 0
$ws.Cells.Item(85, 10).Value = 0  # J85: 32000 -> 0
$ws.Cells.Item(85, 12).Value = 0  # L85: 32000 -> 0
$ws.Cells.Item(85, 14).ClearContents()  # N85: was -34652
# Row 137
$ws.Cells.Item(137, 8).Value = 0  # H137: 42999.668 -> 0
$ws.Cells.Item(137, 10).Value = 0  # J137: 42999.668 -> 0
$ws.Cells.Item(137, 12).Value = 0  # L137: 42999.668 -> 0
$ws.Cells.Item(137, 14).ClearContents()  # N137: was -53199.668
